# Added file read feature
# Populate a new column D on "Лист2" with the additional data read from file,
# matching A1:C3 -> A1:D3 extension.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

$ws.Range("D1").Value = 5
$ws.Range("D2").Value = 7
$ws.Range("D3").Value = 9

$ws.Range("D3").Select() | Out-Null

